$d = $word.ActiveDocument

# Locate the "• +91 7767880235 •" run precisely, then extend into the
# single-space run that immediately follows it. Re-writing just that
# trailing space (while anchoring the edit on the character that follows
# the bullet/phone text) causes the host to coalesce the two identically
# formatted runs into one run whose text is
# "• +91 7767880235 • " (trailing space kept), matching the target diff,
# without touching the separate leading-space run that precedes the
# phone-number run.
$r = $d.Content
$found = $r.Find.Execute("• +91 7767880235 •", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $tail = $d.Range($r.End, $r.End + 1)
    if ($tail.Text -eq " ") {
        $tail.Text = " "
    }
}
